$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new row of data (row 9) between the existing entry for 2012-08-01 (row 8)
# and the blank spacer rows that follow.
$ws.Cells.Item(8, 1).Copy()
$ws.Cells.Item(9, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(9, 1).Value = 41123

$ws.Cells.Item(9, 2).Value = "Implemented Bealto ParallelBitonicSortA, ParallelBitonicSortB2 (first algorithm running faster than CPU Quicksort)"

# Update the active selection to reflect where the author left off editing.
$ws.Activate()
$ws.Range("B14").Select()
